$d = $word.ActiveDocument

$pkgHeader = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:cs="Courier New" w:eastAsia="Courier New" w:hAnsi="Courier New"/><w:color w:val="7f6000"/><w:sz w:val="18"/><w:szCs w:val="18"/><w:rtl w:val="0"/></w:rPr><w:t xml:space="preserve">'
$pkgFooter = '</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

for ($i = 1; $i -le 6; $i++) {
    $old = "<id>p020v_$i</id>"
    $newVisible = "<id>p020v_$i</id>"

    # Locate the run(s) containing the old "<id>...</id>" text (currently
    # split across 3 separate runs: "<id>", "p020v_N", "</id>").
    $rng = $d.Content
    $rng.Find.Execute($old, $true, $false, $false, $false, $false,
                       $true, 1, $false, "", 0)

    $startPos = $rng.Start
    $oldLen = $rng.End - $rng.Start

    # Insert the replacement as a single run (with the formatting of the
    # original first run) at the collapsed start of the match, which
    # leaves the rest of the paragraph (including the trailing run(s))
    # untouched and in place.
    $insertPoint = $d.Range($startPos, $startPos)
    $fragment = $pkgHeader + "&lt;id&gt;p020v_$i&lt;/id&gt;" + $pkgFooter
    $insertPoint.InsertXML($fragment)

    # Now remove the original 3 runs' text, which has been shifted later
    # in the document by the length of the newly inserted text.
    $newLen = $newVisible.Length
    $oldRange = $d.Range($startPos + $newLen, $startPos + $newLen + $oldLen)
    $oldRange.Text = ""
}
